$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 29.41996
$ws.Range("H2").Value = 88.25988000000001
$ws.Range("I2").Value = 0.6207199949605289
$ws.Range("J2").Value = 0.66829493802317
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.107177
$ws.Range("N2").Value = 0.321531
$ws.Range("O2").Value = 0.003526763356587491
$ws.Range("P2").Value = 0.003549676734010809
$ws.Range("Q2").Value = 3.15314305292
$ws.Range("R2").Value = 28.37828747628
$ws.Range("S2").Value = 0.002189132532927965
$ws.Range("T2").Value = 0.002372230992958042
$ws.Range("G3").Value = 29.41996
$ws.Range("H3").Value = 88.25988000000001
$ws.Range("I3").Value = 0.6207199949605289
$ws.Range("J3").Value = 0.66829493802317
$ws.Range("O3").Value = 0.9757678722356318
$ws.Range("P3").Value = 0.9821074349659524
$ws.Range("Q3").Value = 872.3964089780534
$ws.Range("R3").Value = 7851.56768080248
$ws.Range("S3").Value = 0.6056786287367474
$ws.Range("T3").Value = 0.6563374273826657
$ws.Range("D4").Value = "Inflammatory-Mac"
$ws.Range("G4").Value = 29.41996
$ws.Range("H4").Value = 88.25988000000001
$ws.Range("I4").Value = 0.6207199949605289
$ws.Range("J4").Value = 0.66829493802317
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.04072766666666667
$ws.Range("N4").Value = 0.122183
$ws.Range("O4").Value = 0.001340183457265176
$ws.Range("P4").Value = 0.001348890627627329
$ws.Range("Q4").Value = 1.198206324226667
$ws.Range("R4").Value = 10.78385691804
$ws.Range("S4").Value = 0.0008318786688398243
$ws.Range("T4").Value = 0.0009014567783902405
$ws.Range("D5").Value = "MuSCs"
$ws.Range("G5").Value = 29.41996
$ws.Range("H5").Value = 88.25988000000001
$ws.Range("I5").Value = 0.6207199949605289
$ws.Range("J5").Value = 0.66829493802317
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.5885005
$ws.Range("N5").Value = 1.177001
$ws.Range("O5").Value = 0.01936518095051565
$ws.Range("P5").Value = 0.01299399767240936
$ws.Range("Q5").Value = 17.31366116998
$ws.Range("R5").Value = 103.88196701988
$ws.Range("S5").Value = 0.01202035502201381
$ws.Range("T5").Value = 0.008683822869156031
$ws.Range("I6").Value = 0.1515698101047853
$ws.Range("J6").Value = 0.1631868437822795
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.107177
$ws.Range("N6").Value = 0.321531
$ws.Range("O6").Value = 0.003526763356587491
$ws.Range("P6").Value = 0.003549676734010809
$ws.Range("Q6").Value = 0.7699466710343335
$ws.Range("R6").Value = 6.929520039309001
$ws.Range("S6").Value = 0.0005345508522424811
$ws.Range("T6").Value = 0.0005792605426706139
$ws.Range("I7").Value = 0.1515698101047853
$ws.Range("J7").Value = 0.1631868437822795
$ws.Range("O7").Value = 0.9757678722356318
$ws.Range("P7").Value = 0.9821074349659524
$ws.Range("S7").Value = 0.1478969511011051
$ws.Range("T7").Value = 0.1602670125672041
$ws.Range("D8").Value = "Inflammatory-Mac"
$ws.Range("I8").Value = 0.1515698101047853
$ws.Range("J8").Value = 0.1631868437822795
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.04072766666666667
$ws.Range("N8").Value = 0.122183
$ws.Range("O8").Value = 0.001340183457265176
$ws.Range("P8").Value = 0.001348890627627329
$ws.Range("Q8").Value = 0.2925826564374445
$ws.Range("R8").Value = 2.633243907937
$ws.Range("S8").Value = 0.0002031313521232574
$ws.Range("T8").Value = 0.0002201212041300018
$ws.Range("D9").Value = "MuSCs"
$ws.Range("I9").Value = 0.1515698101047853
$ws.Range("J9").Value = 0.1631868437822795
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.5885005
$ws.Range("N9").Value = 1.177001
$ws.Range("O9").Value = 0.01936518095051565
$ws.Range("P9").Value = 0.01299399767240936
$ws.Range("Q9").Value = 4.227716775773167
$ws.Range("R9").Value = 25.366300654639
$ws.Range("S9").Value = 0.002935176799314463
$ws.Range("T9").Value = 0.00212044946827477
$ws.Range("G10").Value = 0.3873096666666667
$ws.Range("H10").Value = 1.161929
$ws.Range("I10").Value = 0.008171692087327698
$ws.Range("J10").Value = 0.008798009571759262
$ws.Range("K10").Value = 2
$ws.Range("L10").Value = 0.6666666666666666
$ws.Range("M10").Value = 0.107177
$ws.Range("N10").Value = 0.321531
$ws.Range("O10").Value = 0.003526763356587491
$ws.Range("P10").Value = 0.003549676734010809
$ws.Range("Q10").Value = 0.04151068814433333
$ws.Range("R10").Value = 0.373596193299
$ws.Range("S10").Value = [double]"2.881962421490327E-05"
$ws.Range("T10").Value = [double]"3.123008988247825E-05"
$ws.Range("G11").Value = 0.3873096666666667
$ws.Range("H11").Value = 1.161929
$ws.Range("I11").Value = 0.008171692087327698
$ws.Range("J11").Value = 0.008798009571759262
$ws.Range("O11").Value = 0.9757678722356318
$ws.Range("P11").Value = 0.9821074349659524
$ws.Range("Q11").Value = 11.48497694634822
$ws.Range("R11").Value = 103.364792517134
$ws.Range("S11").Value = 0.007973674600616498
$ws.Range("T11").Value = 0.008640590613326386
$ws.Range("D12").Value = "Inflammatory-Mac"
$ws.Range("G12").Value = 0.3873096666666667
$ws.Range("H12").Value = 1.161929
$ws.Range("I12").Value = 0.008171692087327698
$ws.Range("J12").Value = 0.008798009571759262
$ws.Range("L12").Value = 0.6666666666666666
$ws.Range("M12").Value = 0.04072766666666667
$ws.Range("N12").Value = 0.122183
$ws.Range("O12").Value = 0.001340183457265176
$ws.Range("P12").Value = 0.001348890627627329
$ws.Range("Q12").Value = 0.01577421900077778
$ws.Range("R12").Value = 0.141967971007
$ws.Range("S12").Value = [double]"1.095156655330132E-05"
$ws.Range("T12").Value = [double]"1.18675526531216E-05"
$ws.Range("D13").Value = "MuSCs"
$ws.Range("G13").Value = 0.3873096666666667
$ws.Range("H13").Value = 1.161929
$ws.Range("I13").Value = 0.008171692087327698
$ws.Range("J13").Value = 0.008798009571759262
$ws.Range("K13").Value = 2
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 0.5885005
$ws.Range("N13").Value = 1.177001
$ws.Range("O13").Value = 0.01936518095051565
$ws.Range("P13").Value = 0.01299399767240936
$ws.Range("Q13").Value = 0.2279319324881666
$ws.Range("R13").Value = 1.367591594929
$ws.Range("S13").Value = 0.0001582462959429979
$ws.Range("T13").Value = 0.0001143213158972752
$ws.Range("G14").Value = 10.122265
$ws.Range("H14").Value = 20.24453
$ws.Range("I14").Value = 0.2135656295858028
$ws.Range("J14").Value = 0.153289545846405
$ws.Range("K14").Value = 2
$ws.Range("L14").Value = 0.6666666666666666
$ws.Range("M14").Value = 0.107177
$ws.Range("N14").Value = 0.321531
$ws.Range("O14").Value = 0.003526763356587491
$ws.Range("P14").Value = 0.003549676734010809
$ws.Range("Q14").Value = 1.084873995905
$ws.Range("R14").Value = 6.50924397543
$ws.Range("S14").Value = 0.0007531954366497468
$ws.Range("T14").Value = 0.000544128334458067
$ws.Range("G15").Value = 10.122265
$ws.Range("H15").Value = 20.24453
$ws.Range("I15").Value = 0.2135656295858028
$ws.Range("J15").Value = 0.153289545846405
$ws.Range("O15").Value = 0.9757678722356318
$ws.Range("P15").Value = 0.9821074349659524
$ws.Range("Q15").Value = 300.1577037060633
$ws.Range("R15").Value = 1800.94622223638
$ws.Range("S15").Value = 0.2083904799636019
$ws.Range("T15").Value = 0.1505468026783086
$ws.Range("D16").Value = "Inflammatory-Mac"
$ws.Range("G16").Value = 10.122265
$ws.Range("H16").Value = 20.24453
$ws.Range("I16").Value = 0.2135656295858028
$ws.Range("J16").Value = 0.153289545846405
$ws.Range("L16").Value = 0.6666666666666666
$ws.Range("M16").Value = 0.04072766666666667
$ws.Range("N16").Value = 0.122183
$ws.Range("O16").Value = 0.001340183457265176
$ws.Range("P16").Value = 0.001348890627627329
$ws.Range("Q16").Value = 0.4122562348316666
$ws.Range("R16").Value = 2.47353740899
$ws.Range("S16").Value = 0.0002862171238113153
$ws.Range("T16").Value = 0.0002067708317054654
$ws.Range("D17").Value = "MuSCs"
$ws.Range("G17").Value = 10.122265
$ws.Range("H17").Value = 20.24453
$ws.Range("I17").Value = 0.2135656295858028
$ws.Range("J17").Value = 0.153289545846405
$ws.Range("K17").Value = 2
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 0.5885005
$ws.Range("N17").Value = 1.177001
$ws.Range("O17").Value = 0.01936518095051565
$ws.Range("P17").Value = 0.01299399767240936
$ws.Range("Q17").Value = 5.956958013632499
$ws.Range("R17").Value = 23.82783205453
$ws.Range("S17").Value = 0.004135737061739871
$ws.Range("T17").Value = 0.001991844001932875
$ws.Range("G18").Value = 0.2830933333333334
$ws.Range("H18").Value = 0.84928
$ws.Range("I18").Value = 0.005972873261555284
$ws.Range("J18").Value = 0.006430662776386256
$ws.Range("K18").Value = 2
$ws.Range("L18").Value = 0.6666666666666666
$ws.Range("M18").Value = 0.107177
$ws.Range("N18").Value = 0.321531
$ws.Range("O18").Value = 0.003526763356587491
$ws.Range("P18").Value = 0.003549676734010809
$ws.Range("Q18").Value = 0.03034109418666667
$ws.Range("R18").Value = 0.27306984768
$ws.Range("S18").Value = [double]"2.106491055239439E-05"
$ws.Range("T18").Value = [double]"2.282677404160764E-05"
$ws.Range("G19").Value = 0.2830933333333334
$ws.Range("H19").Value = 0.84928
$ws.Range("I19").Value = 0.005972873261555284
$ws.Range("J19").Value = 0.006430662776386256
$ws.Range("O19").Value = 0.9757678722356318
$ws.Range("P19").Value = 0.9821074349659524
$ws.Range("Q19").Value = 8.394627572764445
$ws.Range("R19").Value = 75.55164815488
$ws.Range("S19").Value = 0.005828137833560897
$ws.Range("T19").Value = 0.006315601724447736
$ws.Range("D20").Value = "Inflammatory-Mac"
$ws.Range("G20").Value = 0.2830933333333334
$ws.Range("H20").Value = 0.84928
$ws.Range("I20").Value = 0.005972873261555284
$ws.Range("J20").Value = 0.006430662776386256
$ws.Range("L20").Value = 0.6666666666666666
$ws.Range("M20").Value = 0.04072766666666667
$ws.Range("N20").Value = 0.122183
$ws.Range("O20").Value = 0.001340183457265176
$ws.Range("P20").Value = 0.001348890627627329
$ws.Range("Q20").Value = 0.01152973091555556
$ws.Range("R20").Value = 0.10376757824
$ws.Range("S20").Value = [double]"8.00474593747789E-06"
$ws.Range("T20").Value = [double]"8.674260748499358E-06"
$ws.Range("D21").Value = "MuSCs"
$ws.Range("G21").Value = 0.2830933333333334
$ws.Range("H21").Value = 0.84928
$ws.Range("I21").Value = 0.005972873261555284
$ws.Range("J21").Value = 0.006430662776386256
$ws.Range("K21").Value = 2
$ws.Range("L21").Value = 1
$ws.Range("M21").Value = 0.5885005
$ws.Range("N21").Value = 1.177001
$ws.Range("O21").Value = 0.01936518095051565
$ws.Range("P21").Value = 0.01299399767240936
$ws.Range("Q21").Value = 0.1666005682133334
$ws.Range("R21").Value = 0.99960340928
$ws.Range("S21").Value = 0.0001156657715045147
$ws.Range("T21").Value = [double]"8.356001714841255E-05"
